$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = -526.3209500055157
$ws.Cells.Item(2, 2).Value = 30.46895981468204
$ws.Cells.Item(3, 1).Value = -524.4479470372619
$ws.Cells.Item(3, 2).Value = 30.18261735185105
$ws.Cells.Item(4, 1).Value = -522.5619826905364
$ws.Cells.Item(4, 2).Value = 29.90211143955243
$ws.Cells.Item(5, 1).Value = -520.662931608397
$ws.Cells.Item(5, 2).Value = 29.62708459867983
$ws.Cells.Item(6, 1).Value = -518.7507575836784
$ws.Cells.Item(6, 2).Value = 29.35723248935262
$ws.Cells.Item(7, 1).Value = -516.8253772395135
$ws.Cells.Item(7, 2).Value = 29.09225425243264
$ws.Cells.Item(8, 1).Value = -514.8867386148538
$ws.Cells.Item(8, 2).Value = 28.83187799249755
$ws.Cells.Item(9, 1).Value = -512.9348661186879
$ws.Cells.Item(9, 2).Value = 28.57587363536126
$ws.Cells.Item(10, 1).Value = -510.96971159504
$ws.Cells.Item(10, 2).Value = 28.32400164845155
$ws.Cells.Item(11, 1).Value = -508.9913052284678
$ws.Cells.Item(11, 2).Value = 28.07606236471085
$ws.Cells.Item(12, 1).Value = -506.9996563762015
$ws.Cells.Item(12, 2).Value = 27.8318605117626
$ws.Cells.Item(13, 1).Value = -504.9947869090348
$ws.Cells.Item(13, 2).Value = 27.59121681079435
$ws.Cells.Item(14, 1).Value = -502.9767566420569
$ws.Cells.Item(14, 2).Value = 27.35397422070858
$ws.Cells.Item(15, 1).Value = -500.9455843324719
$ws.Cells.Item(15, 2).Value = 27.11997176395216
$ws.Cells.Item(16, 1).Value = -498.9013113393704
$ws.Cells.Item(16, 2).Value = 26.88906441690665
$ws.Cells.Item(17, 1).Value = -496.8440001712912
$ws.Cells.Item(17, 2).Value = 26.66112225137456
$ws.Cells.Item(18, 1).Value = -494.7737074712043
$ws.Cells.Item(18, 2).Value = 26.43601978107654
$ws.Cells.Item(19, 1).Value = -492.6904892925363
$ws.Cells.Item(19, 2).Value = 26.21363952092852
$ws.Cells.Item(20, 1).Value = -490.5944085875326
$ws.Cells.Item(20, 2).Value = 25.99387076134603
$ws.Cells.Item(21, 1).Value = -488.485523249671
$ws.Cells.Item(21, 2).Value = 25.77660873092152
$ws.Cells.Item(22, 1).Value = -486.3639000992166
$ws.Cells.Item(22, 2).Value = 25.5617560512568
$ws.Cells.Item(23, 1).Value = -484.2296142165903
$ws.Cells.Item(23, 2).Value = 25.34922263690954
$ws.Cells.Item(24, 1).Value = -482.0827291540471
$ws.Cells.Item(24, 2).Value = 25.13891955832093
$ws.Cells.Item(25, 1).Value = -479.923318309864
$ws.Cells.Item(25, 2).Value = 24.93076638269421
$ws.Cells.Item(26, 1).Value = -477.7514408685112
$ws.Cells.Item(26, 2).Value = 24.72468144344471
$ws.Cells.Item(27, 1).Value = -475.5671951964998
$ws.Cells.Item(27, 2).Value = 24.52059835085536
$ws.Cells.Item(28, 1).Value = -473.3706274944925
$ws.Cells.Item(28, 2).Value = 24.31844065195307
$ws.Cells.Item(29, 1).Value = -471.1618247999816
$ws.Cells.Item(29, 2).Value = 24.11814501345723
$ws.Cells.Item(30, 1).Value = -468.940872925242
$ws.Cells.Item(30, 2).Value = 23.91965211721636
$ws.Cells.Item(31, 1).Value = -466.7078104848711
$ws.Cells.Item(31, 2).Value = 23.72289229025026
$ws.Cells.Item(32, 1).Value = -464.462736640075
$ws.Cells.Item(32, 2).Value = 23.52781537446108
$ws.Cells.Item(33, 1).Value = -462.2057154767515
$ws.Cells.Item(33, 2).Value = 23.33436428155985
$ws.Cells.Item(34, 1).Value = -459.9368368608303
$ws.Cells.Item(34, 2).Value = 23.1424909539227
$ws.Cells.Item(35, 1).Value = -457.6561677880746
$ws.Cells.Item(35, 2).Value = 22.95214397455752
$ws.Cells.Item(36, 1).Value = -455.3637651849463
$ws.Cells.Item(36, 2).Value = 22.76327191687544
$ws.Cells.Item(37, 1).Value = -453.0597152692862
$ws.Cells.Item(37, 2).Value = 22.57583236588632
$ws.Cells.Item(38, 1).Value = -450.7441062658509
$ws.Cells.Item(38, 2).Value = 22.38978583332456
$ws.Cells.Item(39, 1).Value = -448.4169828563089
$ws.Cells.Item(39, 2).Value = 22.20508311565503
$ws.Cells.Item(40, 1).Value = -446.0784483024984
$ws.Cells.Item(40, 2).Value = 22.02169219641556
$ws.Cells.Item(41, 1).Value = -443.7285680464067
$ws.Cells.Item(41, 2).Value = 21.83957295414289
$ws.Cells.Item(42, 1).Value = -441.3673951192543
$ws.Cells.Item(42, 2).Value = 21.65868364459233
$ws.Cells.Item(43, 1).Value = -438.995019040299
$ws.Cells.Item(43, 2).Value = 21.47899350834041
$ws.Cells.Item(44, 1).Value = -436.6115107136291
$ws.Cells.Item(44, 2).Value = 21.3004682654454
$ws.Cells.Item(45, 1).Value = -434.216934671644
$ws.Cells.Item(45, 2).Value = 21.12307336968293
$ws.Cells.Item(46, 1).Value = -431.8113915930381
$ws.Cells.Item(46, 2).Value = 20.94678469772343
$ws.Cells.Item(47, 1).Value = -429.3949105346652
$ws.Cells.Item(47, 2).Value = 20.77156197876188
$ws.Cells.Item(48, 1).Value = -426.9675926260201
$ws.Cells.Item(48, 2).Value = 20.59738293802198
$ws.Cells.Item(49, 1).Value = -424.5295056155436
$ws.Cells.Item(49, 2).Value = 20.42421933096963
$ws.Cells.Item(50, 1).Value = -422.0806972971806
$ws.Cells.Item(50, 2).Value = 20.25203834685896
$ws.Cells.Item(51, 1).Value = -419.6212692187369
$ws.Cells.Item(51, 2).Value = 20.08082152047011
$ws.Cells.Item(52, 1).Value = -417.151273553826
$ws.Cells.Item(52, 2).Value = 19.91053956337895
$ws.Cells.Item(53, 1).Value = -414.6708024864402
$ws.Cells.Item(53, 2).Value = 19.74117324719001
$ws.Cells.Item(54, 1).Value = -412.1798986416407
$ws.Cells.Item(54, 2).Value = 19.57269278185451
$ws.Cells.Item(55, 1).Value = -409.6786432714717
$ws.Cells.Item(55, 2).Value = 19.4050785638129
$ws.Cells.Item(56, 1).Value = -407.1671067576169
$ws.Cells.Item(56, 2).Value = 19.23830850141663
$ws.Cells.Item(57, 1).Value = -404.6453569094813
$ws.Cells.Item(57, 2).Value = 19.07236163077455
$ws.Cells.Item(58, 1).Value = -402.1134697778352
$ws.Cells.Item(58, 2).Value = 18.90721853039033
$ws.Cells.Item(59, 1).Value = -399.5715032431371
$ws.Cells.Item(59, 2).Value = 18.7428574633572
$ws.Cells.Item(60, 1).Value = -397.0195356616037
$ws.Cells.Item(60, 2).Value = 18.57926092960198
$ws.Cells.Item(61, 1).Value = -394.4576340006365
$ws.Cells.Item(61, 2).Value = 18.416410285178
$ws.Cells.Item(62, 1).Value = -391.8858629321753
$ws.Cells.Item(62, 2).Value = 18.25428658636625
$ws.Cells.Item(63, 1).Value = -389.3042893209173
$ws.Cells.Item(63, 2).Value = 18.09287218552746
$ws.Cells.Item(64, 1).Value = -386.7129923121398
$ws.Cells.Item(64, 2).Value = 17.93215248862646
$ws.Cells.Item(65, 1).Value = -384.1120294816633
$ws.Cells.Item(65, 2).Value = 17.77210876417702
$ws.Cells.Item(66, 1).Value = -381.5014635421405
$ws.Cells.Item(66, 2).Value = 17.61272418891071
$ws.Cells.Item(67, 1).Value = -378.8813853470432
$ws.Cells.Item(67, 2).Value = 17.45398801665599
$ws.Cells.Item(68, 1).Value = -376.2518483921783
$ws.Cells.Item(68, 2).Value = 17.29588223748921
$ws.Cells.Item(69, 1).Value = -373.6129180285574
$ws.Cells.Item(69, 2).Value = 17.13839185755111
$ws.Cells.Item(70, 1).Value = -370.9646569795685
$ws.Cells.Item(70, 2).Value = 16.98150155760036
$ws.Cells.Item(71, 1).Value = -368.307154618633
$ws.Cells.Item(71, 2).Value = 16.82520210954121
$ws.Cells.Item(72, 1).Value = -365.6404627313224
$ws.Cells.Item(72, 2).Value = 16.66947717531342
$ws.Cells.Item(73, 1).Value = -362.9646410439397
$ws.Cells.Item(73, 2).Value = 16.51431179695347
$ws.Cells.Item(74, 1).Value = -360.2797716173147
$ws.Cells.Item(74, 2).Value = 16.35969652565372
$ws.Cells.Item(75, 1).Value = -357.5859128926795
$ws.Cells.Item(75, 2).Value = 16.20561714842691
$ws.Cells.Item(76, 1).Value = -354.8831560323578
$ws.Cells.Item(76, 2).Value = 16.05206654012009
$ws.Cells.Item(77, 1).Value = -352.1715348921624
$ws.Cells.Item(77, 2).Value = 15.89902619758993
$ws.Cells.Item(78, 1).Value = -349.4511317474088
$ws.Cells.Item(78, 2).Value = 15.74648766643109
$ws.Cells.Item(79, 1).Value = -346.7220130370501
$ws.Cells.Item(79, 2).Value = 15.59443975583335
$ws.Cells.Item(80, 1).Value = -343.9842602169075
$ws.Cells.Item(80, 2).Value = 15.44287450541905
$ws.Cells.Item(81, 1).Value = -341.2379157727331
$ws.Cells.Item(81, 2).Value = 15.29177637471872
$ws.Cells.Item(82, 1).Value = -338.4830547944313
$ws.Cells.Item(82, 2).Value = 15.14113670805183
$ws.Cells.Item(83, 1).Value = -335.7197559666628
$ws.Cells.Item(83, 2).Value = 14.99094755830849
$ws.Cells.Item(84, 1).Value = -332.9480846806693
$ws.Cells.Item(84, 2).Value = 14.84119892222046
$ws.Cells.Item(85, 1).Value = -330.1681071475258
$ws.Cells.Item(85, 2).Value = 14.69188091774833
$ws.Cells.Item(86, 1).Value = -327.3798862233397
$ws.Cells.Item(86, 2).Value = 14.5429833138509
$ws.Cells.Item(87, 1).Value = -324.5834879778962
$ws.Cells.Item(87, 2).Value = 14.39449664021209
$ws.Cells.Item(88, 1).Value = -321.7789892566498
$ws.Cells.Item(88, 2).Value = 14.24641402267901
$ws.Cells.Item(89, 1).Value = -318.9664507762616
$ws.Cells.Item(89, 2).Value = 14.09872511810034
$ws.Cells.Item(90, 1).Value = -316.1459347166937
$ws.Cells.Item(90, 2).Value = 13.95142048601543
$ws.Cells.Item(91, 1).Value = -313.3175411126151
$ws.Cells.Item(91, 2).Value = 13.80449833515382
$ws.Cells.Item(92, 1).Value = -310.4812898374654
$ws.Cells.Item(92, 2).Value = 13.65794077094083
$ws.Cells.Item(93, 1).Value = -307.6372920508073
$ws.Cells.Item(93, 2).Value = 13.51174861300277
$ws.Cells.Item(94, 1).Value = -304.7855965737814
$ws.Cells.Item(94, 2).Value = 13.36591005194398
$ws.Cells.Item(95, 1).Value = -301.9262682298449
$ws.Cells.Item(95, 2).Value = 13.22041699459906
$ws.Cells.Item(96, 1).Value = -299.0593834678788
$ws.Cells.Item(96, 2).Value = 13.07526373922262
$ws.Cells.Item(97, 1).Value = -296.1850203820026
$ws.Cells.Item(97, 2).Value = 12.93044456236887
$ws.Cells.Item(98, 1).Value = -293.3032201123488
$ws.Cells.Item(98, 2).Value = 12.78594730028135
$ws.Cells.Item(99, 1).Value = -290.4140890940761
$ws.Cells.Item(99, 2).Value = 12.641772119318
$ws.Cells.Item(100, 1).Value = -287.5176693335586
$ws.Cells.Item(100, 2).Value = 12.49790704912911
$ws.Cells.Item(101, 1).Value = -284.6140372226678
$ws.Cells.Item(101, 2).Value = 12.35434696350926
$ws.Cells.Item(102, 1).Value = -281.7032703667011
$ws.Cells.Item(102, 2).Value = 12.21108673983932
$ws.Cells.Item(103, 1).Value = -278.7854335852587
$ws.Cells.Item(103, 2).Value = 12.06811923078231
$ws.Cells.Item(104, 1).Value = -275.8605865910004
$ws.Cells.Item(104, 2).Value = 11.92543626770419
$ws.Cells.Item(105, 1).Value = -272.9288117756543
$ws.Cells.Item(105, 2).Value = 11.78303418881602
$ws.Cells.Item(106, 1).Value = -269.9901748618988
$ws.Cells.Item(106, 2).Value = 11.6409060992191
$ws.Cells.Item(107, 1).Value = -267.0447582494876
$ws.Cells.Item(107, 2).Value = 11.49904849009149
$ws.Cells.Item(108, 1).Value = -264.0926103014244
$ws.Cells.Item(108, 2).Value = 11.35745141530824
$ws.Cells.Item(109, 1).Value = -261.1338194236249
$ws.Cells.Item(109, 2).Value = 11.21611260405959
$ws.Cells.Item(110, 1).Value = -258.1684566317879
$ws.Cells.Item(110, 2).Value = 11.07502662001487
$ws.Cells.Item(111, 1).Value = -255.1965819322419
$ws.Cells.Item(111, 2).Value = 10.93418594664655
$ws.Cells.Item(112, 1).Value = -252.2182778904401
$ws.Cells.Item(112, 2).Value = 10.79358744966037
$ws.Cells.Item(113, 1).Value = -249.2336105298114
$ws.Cells.Item(113, 2).Value = 10.65322482627247
$ws.Cells.Item(114, 1).Value = -246.2426507075173
$ws.Cells.Item(114, 2).Value = 10.51309296544647
$ws.Cells.Item(115, 1).Value = -243.2454814958002
$ws.Cells.Item(115, 2).Value = 10.37318889458071
$ws.Cells.Item(116, 1).Value = -240.2421635749985
$ws.Cells.Item(116, 2).Value = 10.23350545782563
$ws.Cells.Item(117, 1).Value = -237.2327790037845
$ws.Cells.Item(117, 2).Value = 10.09403985710317
$ws.Cells.Item(118, 1).Value = -234.2173886192221
$ws.Cells.Item(118, 2).Value = 9.954785076102587
$ws.Cells.Item(119, 1).Value = -231.1960753389904
$ws.Cells.Item(119, 2).Value = 9.815738391174136
$ws.Cells.Item(120, 1).Value = -228.1689049313135
$ws.Cells.Item(120, 2).Value = 9.67689400021129
$ws.Cells.Item(121, 1).Value = -225.1359556525492
$ws.Cells.Item(121, 2).Value = 9.53824819691593
$ws.Cells.Item(122, 1).Value = -222.0973040578197
$ws.Cells.Item(122, 2).Value = 9.399797396932096
$ws.Cells.Item(123, 1).Value = -219.053017207747
$ws.Cells.Item(123, 2).Value = 9.261535917649164
$ws.Cells.Item(124, 1).Value = -216.0031779981518
$ws.Cells.Item(124, 2).Value = 9.123461279835389
$ws.Cells.Item(125, 1).Value = -212.9478532963885
$ws.Cells.Item(125, 2).Value = 8.985567917443968
$ws.Cells.Item(126, 1).Value = -209.8871148678043
$ws.Cells.Item(126, 2).Value = 8.847851382484558
$ws.Cells.Item(127, 1).Value = -206.8210523162889
$ws.Cells.Item(127, 2).Value = 8.710310305981931
$ws.Cells.Item(128, 1).Value = -203.7497266213126
$ws.Cells.Item(128, 2).Value = 8.572938259111648
$ws.Cells.Item(129, 1).Value = -200.6732152873514
$ws.Cells.Item(129, 2).Value = 8.435731945032444
$ws.Cells.Item(130, 1).Value = -197.5916029841501
$ws.Cells.Item(130, 2).Value = 8.29868906171057
$ws.Cells.Item(131, 1).Value = -194.5049558060952
$ws.Cells.Item(131, 2).Value = 8.16180435802827
$ws.Cells.Item(132, 1).Value = -191.4133470276406
$ws.Cells.Item(132, 2).Value = 8.025073591715993
$ws.Cells.Item(133, 1).Value = -188.3168605900732
$ws.Cells.Item(133, 2).Value = 7.888494598352859
$ws.Cells.Item(134, 1).Value = -185.2155689854326
$ws.Cells.Item(134, 2).Value = 7.752063232587599
$ws.Cells.Item(135, 1).Value = -182.1095511846597
$ws.Cells.Item(135, 2).Value = 7.615776364851915
$ws.Cells.Item(136, 1).Value = -178.9988856200462
$ws.Cells.Item(136, 2).Value = 7.479630913255447
$ws.Cells.Item(137, 1).Value = -175.8836565570712
$ws.Cells.Item(137, 2).Value = 7.343624816954421
$ws.Cells.Item(138, 1).Value = -172.7639093763262
$ws.Cells.Item(138, 2).Value = 7.20774904790652
$ws.Cells.Item(139, 1).Value = -169.6397563311421
$ws.Cells.Item(139, 2).Value = 7.072006603252816
$ws.Cells.Item(140, 1).Value = -166.5112599039525
$ws.Cells.Item(140, 2).Value = 6.936391504512319
$ws.Cells.Item(141, 1).Value = -163.3785044221037
$ws.Cells.Item(141, 2).Value = 6.800901801912753
$ws.Cells.Item(142, 1).Value = -160.2415637277344
$ws.Cells.Item(142, 2).Value = 6.665533566229775
$ws.Cells.Item(143, 1).Value = -157.1005166147034
$ws.Cells.Item(143, 2).Value = 6.53028390357689
$ws.Cells.Item(144, 1).Value = -153.9554373999102
$ws.Cells.Item(144, 2).Value = 6.395148930852366
$ws.Cells.Item(145, 1).Value = -150.8064161877569
$ws.Cells.Item(145, 2).Value = 6.260127769820247
$ws.Cells.Item(146, 1).Value = -147.653521759509
$ws.Cells.Item(146, 2).Value = 6.125215597965129
$ws.Cells.Item(147, 1).Value = -144.4968275960122
$ws.Cells.Item(147, 2).Value = 5.99040862788886
$ws.Cells.Item(148, 1).Value = -141.3364305742743
$ws.Cells.Item(148, 2).Value = 5.855706992162077
$ws.Cells.Item(149, 1).Value = -138.1723992008402
$ws.Cells.Item(149, 2).Value = 5.721105946811651
$ws.Cells.Item(150, 1).Value = -135.0048183303786
$ws.Cells.Item(150, 2).Value = 5.58660371431647
$ws.Cells.Item(151, 1).Value = -131.83375721047
$ws.Cells.Item(151, 2).Value = 5.452195584545723
$ws.Cells.Item(152, 1).Value = -128.6593016616121
$ws.Cells.Item(152, 2).Value = 5.317879800062957
$ws.Cells.Item(153, 1).Value = -125.4815424550374
$ws.Cells.Item(153, 2).Value = 5.183655596977058
$ws.Cells.Item(154, 1).Value = -122.3005542282172
$ws.Cells.Item(154, 2).Value = 5.04951929548339
$ws.Cells.Item(155, 1).Value = -119.1164060595936
$ws.Cells.Item(155, 2).Value = 4.91546626824067
$ws.Cells.Item(156, 1).Value = -115.9291840836481
$ws.Cells.Item(156, 2).Value = 4.781494811956255
$ws.Cells.Item(157, 1).Value = -112.7389850393786
$ws.Cells.Item(157, 2).Value = 4.64760517357079
$ws.Cells.Item(158, 1).Value = -109.5458674599286
$ws.Cells.Item(158, 2).Value = 4.51379083126209
$ws.Cells.Item(159, 1).Value = -106.3499282390868
$ws.Cells.Item(159, 2).Value = 4.380052053604679
$ws.Cells.Item(160, 1).Value = -103.1512427539577
$ws.Cells.Item(160, 2).Value = 4.246385246415519
$ws.Cells.Item(161, 1).Value = -99.94989136854019
$ws.Cells.Item(161, 2).Value = 4.112787801649064
$ws.Cells.Item(162, 1).Value = -96.74595524573519
$ws.Cells.Item(162, 2).Value = 3.979257114265228
$ws.Cells.Item(163, 1).Value = -93.53952591090258
$ws.Cells.Item(163, 2).Value = 3.845792518928382
$ws.Cells.Item(164, 1).Value = -90.3306734055017
$ws.Cells.Item(164, 2).Value = 3.712389509287564
$ws.Cells.Item(165, 1).Value = -87.11948939389166
$ws.Cells.Item(165, 2).Value = 3.579047435631556
$ws.Cells.Item(166, 1).Value = -83.90604429450177
$ws.Cells.Item(166, 2).Value = 3.445761813140363
$ws.Cells.Item(167, 1).Value = -80.69043012719757
$ws.Cells.Item(167, 2).Value = 3.312532005031693
$ws.Cells.Item(168, 1).Value = -77.47273335161529
$ws.Cells.Item(168, 2).Value = 3.179356420779543
$ws.Cells.Item(169, 1).Value = -74.25302411862549
$ws.Cells.Item(169, 2).Value = 3.046230608174044
$ws.Cells.Item(170, 1).Value = -71.03138946779978
$ws.Cells.Item(170, 2).Value = 2.913152990647172
$ws.Cells.Item(171, 1).Value = -67.80792701483203
$ws.Cells.Item(171, 2).Value = 2.780123906420734
$ws.Cells.Item(172, 1).Value = -64.58269632875601
$ws.Cells.Item(172, 2).Value = 2.647137019076048
$ws.Cells.Item(173, 1).Value = -61.35579513112401
$ws.Cells.Item(173, 2).Value = 2.514192678535426
$ws.Cells.Item(174, 1).Value = -58.12729945723294
$ws.Cells.Item(174, 2).Value = 2.381287424933994
$ws.Cells.Item(175, 1).Value = -54.89730711333954
$ws.Cells.Item(175, 2).Value = 2.248421612557848
$ws.Cells.Item(176, 1).Value = -51.66588347739155
$ws.Cells.Item(176, 2).Value = 2.115589889186911
$ws.Cells.Item(177, 1).Value = -48.43311547483003
$ws.Cells.Item(177, 2).Value = 1.982790716671695
$ws.Cells.Item(178, 1).Value = -45.1990956583507
$ws.Cells.Item(178, 2).Value = 1.850023508759819
$ws.Cells.Item(179, 1).Value = -41.9639058674862
$ws.Cells.Item(179, 2).Value = 1.717285781251364
$ws.Cells.Item(180, 1).Value = -38.72762217724998
$ws.Cells.Item(180, 2).Value = 1.584574106496285
$ws.Cells.Item(181, 1).Value = -35.49033187063916
$ws.Cells.Item(181, 2).Value = 1.451886958721586
$ws.Cells.Item(182, 1).Value = -32.25212762785644
$ws.Cells.Item(182, 2).Value = 1.31922376246698
$ws.Cells.Item(183, 1).Value = -29.01307509127257
$ws.Cells.Item(183, 2).Value = 1.186579204195945
$ws.Cells.Item(184, 1).Value = -25.77328310259912
$ws.Cells.Item(184, 2).Value = 1.053955557116034
$ws.Cells.Item(185, 1).Value = -22.5328117746925
$ws.Cells.Item(185, 2).Value = 0.9213465666531531
$ws.Cells.Item(186, 1).Value = -19.29176483067081
$ws.Cells.Item(186, 2).Value = 0.7887535603924168
$ws.Cells.Item(187, 1).Value = -16.0502022168778
$ws.Cells.Item(187, 2).Value = 0.6561702905680197
$ws.Cells.Item(188, 1).Value = -12.80822768719755
$ws.Cells.Item(188, 2).Value = 0.5235980868076636
$ws.Cells.Item(189, 1).Value = -9.565934266600436
$ws.Cells.Item(189, 2).Value = 0.3910363847762161
$ws.Cells.Item(190, 1).Value = -6.323376575789477
$ws.Cells.Item(190, 2).Value = 0.2584779959552746
$ws.Cells.Item(191, 1).Value = -3.080658186022802
$ws.Cells.Item(191, 2).Value = 0.1259242514358974
$ws.Cells.Item(192, 1).Value = 0.1621387999721795
$ws.Cells.Item(192, 2).Value = -0.006627303440934691
$ws.Cells.Item(193, 1).Value = 3.404932694777011
$ws.Cells.Item(193, 2).Value = -0.1391791229783856
$ws.Cells.Item(194, 1).Value = 6.64764139343102
$ws.Cells.Item(194, 2).Value = -0.2717336615223958
$ws.Cells.Item(195, 1).Value = 9.890166270743846
$ws.Cells.Item(195, 2).Value = -0.4042905343660159
$ws.Cells.Item(196, 1).Value = 13.13243680176767
$ws.Cells.Item(196, 2).Value = -0.5368540892902213
$ws.Cells.Item(197, 1).Value = 16.37437626644921
$ws.Cells.Item(197, 2).Value = -0.6694277292077621
$ws.Cells.Item(198, 1).Value = 19.61588648083842
$ws.Cells.Item(198, 2).Value = -0.8020110712669887
$ws.Cells.Item(199, 1).Value = 22.85689095781071
$ws.Cells.Item(199, 2).Value = -0.9346075207876816
$ws.Cells.Item(200, 1).Value = 26.09729111690158
$ws.Cells.Item(200, 2).Value = -1.067216696367124
$ws.Cells.Item(201, 1).Value = 29.33702208733256
$ws.Cells.Item(201, 2).Value = -1.199843901295098
$ws.Cells.Item(202, 1).Value = 32.57598511398317
$ws.Cells.Item(202, 2).Value = -1.332488757457292
$ws.Cells.Item(203, 1).Value = 35.8141095834069
$ws.Cells.Item(203, 2).Value = -1.465155627039892
$ws.Cells.Item(204, 1).Value = 39.05130279237452
$ws.Cells.Item(204, 2).Value = -1.59784508329234
$ws.Cells.Item(205, 1).Value = 42.28747721807694
$ws.Cells.Item(205, 2).Value = -1.730558649713433
$ws.Cells.Item(206, 1).Value = 45.52255116469627
$ws.Cells.Item(206, 2).Value = -1.863298801021858
$ws.Cells.Item(207, 1).Value = 48.75644877760389
$ws.Cells.Item(207, 2).Value = -1.996068965163015
$ws.Cells.Item(208, 1).Value = 51.98908800680759
$ws.Cells.Item(208, 2).Value = -2.128871627744429
$ws.Cells.Item(209, 1).Value = 55.22036555414184
$ws.Cells.Item(209, 2).Value = -2.261705472706421
$ws.Cells.Item(210, 1).Value = 58.45021599682435
$ws.Cells.Item(210, 2).Value = -2.394575844518177
$ws.Cells.Item(211, 1).Value = 61.67855783752847
$ws.Cells.Item(211, 2).Value = -2.527485241859197
$ws.Cells.Item(212, 1).Value = 64.90529887212591
$ws.Cells.Item(212, 2).Value = -2.660434262078227
$ws.Cells.Item(213, 1).Value = 68.13035187082468
$ws.Cells.Item(213, 2).Value = -2.793424459675225
$ws.Cells.Item(214, 1).Value = 71.3536410755422
$ws.Cells.Item(214, 2).Value = -2.926459299803378
$ws.Cells.Item(215, 1).Value = 74.57508495918665
$ws.Cells.Item(215, 2).Value = -3.0595413038509
$ws.Cells.Item(216, 1).Value = 77.79459689246265
$ws.Cells.Item(216, 2).Value = -3.192672042490938
$ws.Cells.Item(217, 1).Value = 81.01209007147085
$ws.Cells.Item(217, 2).Value = -3.325853092126767
$ws.Cells.Item(218, 1).Value = 84.22749412617605
$ws.Cells.Item(218, 2).Value = -3.459088908228548
$ws.Cells.Item(219, 1).Value = 87.44071160040333
$ws.Cells.Item(219, 2).Value = -3.592379165534927
$ws.Cells.Item(220, 1).Value = 90.65166686290817
$ws.Cells.Item(220, 2).Value = -3.725727377003012
$ws.Cells.Item(221, 1).Value = 93.86027896720056
$ws.Cells.Item(221, 2).Value = -3.859136106770573
$ws.Cells.Item(222, 1).Value = 97.06646117034272
$ws.Cells.Item(222, 2).Value = -3.992606970356985
$ws.Cells.Item(223, 1).Value = 100.2701385231302
$ws.Cells.Item(223, 2).Value = -4.126143507617892
$ws.Cells.Item(224, 1).Value = 103.4712300518529
$ws.Cells.Item(224, 2).Value = -4.259748314968111
$ws.Cells.Item(225, 1).Value = 106.6696438182927
$ws.Cells.Item(225, 2).Value = -4.393422070398856
$ws.Cells.Item(226, 1).Value = 109.8653102535903
$ws.Cells.Item(226, 2).Value = -4.527169315912123
$ws.Cells.Item(227, 1).Value = 113.0581378027537
$ws.Cells.Item(227, 2).Value = -4.660990745154594
$ws.Cells.Item(228, 1).Value = 116.2480460176158
$ws.Cells.Item(228, 2).Value = -4.794888991194417
$ws.Cells.Item(229, 1).Value = 119.4349650796665
$ws.Cells.Item(229, 2).Value = -4.928868646478103
$ws.Cells.Item(230, 1).Value = 122.6187928440397
$ws.Cells.Item(230, 2).Value = -5.062928489677017
$ws.Cells.Item(231, 1).Value = 125.7994709611893
$ws.Cells.Item(231, 2).Value = -5.197075078963247
$ws.Cells.Item(232, 1).Value = 128.9769025829471
$ws.Cells.Item(232, 2).Value = -5.331308189576799
$ws.Cells.Item(233, 1).Value = 132.1510135880785
$ws.Cells.Item(233, 2).Value = -5.465631486761335
$ws.Cells.Item(234, 1).Value = 135.3217232357661
$ws.Cells.Item(234, 2).Value = -5.60004769993192
$ws.Cells.Item(235, 1).Value = 138.4889630056848
$ws.Cells.Item(235, 2).Value = -5.73456151275823
$ws.Cells.Item(236, 1).Value = 141.6526360831735
$ws.Cells.Item(236, 2).Value = -5.869172750512213
$ws.Cells.Item(237, 1).Value = 144.8126687166878
$ws.Cells.Item(237, 2).Value = -6.003885146566021
$ws.Cells.Item(238, 1).Value = 147.9689804807463
$ws.Cells.Item(238, 2).Value = -6.138701499495634
$ws.Cells.Item(239, 1).Value = 151.1214921216741
$ws.Cells.Item(239, 2).Value = -6.273624603925627
$ws.Cells.Item(240, 1).Value = 154.2701296314933
$ws.Cells.Item(240, 2).Value = -6.408658267612298
$ws.Cells.Item(241, 1).Value = 157.4148022616484
$ws.Cells.Item(241, 2).Value = -6.543803363536482
$ws.Cells.Item(242, 1).Value = 160.5554475998287
$ws.Cells.Item(242, 2).Value = -6.67906571017074
$ws.Cells.Item(243, 1).Value = 163.6919750715018
$ws.Cells.Item(243, 2).Value = -6.814446219283036
$ws.Cells.Item(244, 1).Value = 166.8242995868939
$ws.Cells.Item(244, 2).Value = -6.949946807418678
$ws.Cells.Item(245, 1).Value = 169.9523756377142
$ws.Cells.Item(245, 2).Value = -7.085576363375359
$ws.Cells.Item(246, 1).Value = 173.0760908652595
$ws.Cells.Item(246, 2).Value = -7.221331874583646
$ws.Cells.Item(247, 1).Value = 176.1953773338094
$ws.Cells.Item(247, 2).Value = -7.357218308949079
$ws.Cells.Item(248, 1).Value = 179.3101617788743
$ws.Cells.Item(248, 2).Value = -7.493239669005756
$ws.Cells.Item(249, 1).Value = 182.4203652876251
$ws.Cells.Item(249, 2).Value = -7.629398991589696
$ws.Cells.Item(250, 1).Value = 185.5259150960957
$ws.Cells.Item(250, 2).Value = -7.765700331333
$ws.Cells.Item(251, 1).Value = 188.6267209508647
$ws.Cells.Item(251, 2).Value = -7.902144778203033
$ws.Cells.Item(252, 1).Value = 191.7227269364649
$ws.Cells.Item(252, 2).Value = -8.038739474350313
$ws.Cells.Item(253, 1).Value = 194.8138378216867
$ws.Cells.Item(253, 2).Value = -8.175484543871034
$ws.Cells.Item(254, 1).Value = 197.8999863518384
$ws.Cells.Item(254, 2).Value = -8.312385192329828
$ws.Cells.Item(255, 1).Value = 200.9810944578619
$ws.Cells.Item(255, 2).Value = -8.449444629006658
$ws.Cells.Item(256, 1).Value = 204.0570892974833
$ws.Cells.Item(256, 2).Value = -8.586667131342351
$ws.Cells.Item(257, 1).Value = 207.1278926451112
$ws.Cells.Item(257, 2).Value = -8.72405599227374
$ws.Cells.Item(258, 1).Value = 210.1934266821121
$ws.Cells.Item(258, 2).Value = -8.861614526100606
$ws.Cells.Item(259, 1).Value = 213.2536304406424
$ws.Cells.Item(259, 2).Value = -8.999349168460569
$ws.Cells.Item(260, 1).Value = 216.3084143254626
$ws.Cells.Item(260, 2).Value = -9.137261294155499
$ws.Cells.Item(261, 1).Value = 219.357706582019
$ws.Cells.Item(261, 2).Value = -9.275355357858135
$ws.Cells.Item(262, 1).Value = 222.4014293796985
$ws.Cells.Item(262, 2).Value = -9.413634848090451
$ws.Cells.Item(263, 1).Value = 225.4395159321172
$ws.Cells.Item(263, 2).Value = -9.552105384209032
$ws.Cells.Item(264, 1).Value = 228.4718837673471
$ws.Cells.Item(264, 2).Value = -9.690769474179636
$ws.Cells.Item(265, 1).Value = 231.4984659070298
$ws.Cells.Item(265, 2).Value = -9.829632856103713
$ws.Cells.Item(266, 1).Value = 234.5191906639038
$ws.Cells.Item(266, 2).Value = -9.968700250902454
$ws.Cells.Item(267, 1).Value = 237.5339864325836
$ws.Cells.Item(267, 2).Value = -10.10797643607135
$ws.Cells.Item(268, 1).Value = 240.5427706694153
$ws.Cells.Item(268, 2).Value = -10.24746412482733
$ws.Cells.Item(269, 1).Value = 243.5454766767629
$ws.Cells.Item(269, 2).Value = -10.38716929463351
$ws.Cells.Item(270, 1).Value = 246.5420330488399
$ws.Cells.Item(270, 2).Value = -10.52709689637324
$ws.Cells.Item(271, 1).Value = 249.5323684509479
$ws.Cells.Item(271, 2).Value = -10.66725194471096
$ws.Cells.Item(272, 1).Value = 252.5163942853127
$ws.Cells.Item(272, 2).Value = -10.80763633503746
$ws.Cells.Item(273, 1).Value = 255.4940621730374
$ws.Cells.Item(273, 2).Value = -10.94825946015381
$ws.Cells.Item(274, 1).Value = 258.4652887660231
$ws.Cells.Item(274, 2).Value = -11.08912444819188
$ws.Cells.Item(275, 1).Value = 261.4299977916849
$ws.Cells.Item(275, 2).Value = -11.23023547534267
$ws.Cells.Item(276, 1).Value = 264.3881173961906
$ws.Cells.Item(276, 2).Value = -11.37159793229135
$ws.Cells.Item(277, 1).Value = 267.3395995102687
$ws.Cells.Item(277, 2).Value = -11.5132215732819
$ws.Cells.Item(278, 1).Value = 270.2843442681515
$ws.Cells.Item(278, 2).Value = -11.65510652834324
$ws.Cells.Item(279, 1).Value = 273.222291805249
$ws.Cells.Item(279, 2).Value = -11.79726057696618
$ws.Cells.Item(280, 1).Value = 276.1533827332928
$ws.Cells.Item(280, 2).Value = -11.93969158933592
$ws.Cells.Item(281, 1).Value = 279.0775227392694
$ws.Cells.Item(281, 2).Value = -12.08240099627977
$ws.Cells.Item(282, 1).Value = 281.9946640224617
$ws.Cells.Item(282, 2).Value = -12.22539904102333
$ws.Cells.Item(283, 1).Value = 284.9047293917728
$ws.Cells.Item(283, 2).Value = -12.36869063824693
$ws.Cells.Item(284, 1).Value = 287.8076364689289
$ws.Cells.Item(284, 2).Value = -12.51227963562049
$ws.Cells.Item(285, 1).Value = 290.7033315901436
$ws.Cells.Item(285, 2).Value = -12.65617552725074
$ws.Cells.Item(286, 1).Value = 293.5917427282365
$ws.Cells.Item(286, 2).Value = -12.80038468327549
$ws.Cells.Item(287, 1).Value = 296.4727827333097
$ws.Cells.Item(287, 2).Value = -12.94491005521942
$ws.Cells.Item(288, 1).Value = 299.3464091929567
$ws.Cells.Item(288, 2).Value = -13.0897637777694
$ws.Cells.Item(289, 1).Value = 302.2125388641067
$ws.Cells.Item(289, 2).Value = -13.23495029437471
$ws.Cells.Item(290, 1).Value = 305.0711069627077
$ws.Cells.Item(290, 2).Value = -13.38047745804979
$ws.Cells.Item(291, 1).Value = 307.922041805776
$ws.Cells.Item(291, 2).Value = -13.52635220529397
$ws.Cells.Item(292, 1).Value = 310.7652669118792
$ws.Cells.Item(292, 2).Value = -13.67258037288721
$ws.Cells.Item(293, 1).Value = 313.6007282476484
$ws.Cells.Item(293, 2).Value = -13.81917257341002
$ws.Cells.Item(294, 1).Value = 316.4283544733159
$ws.Cells.Item(294, 2).Value = -13.96613613677125
$ws.Cells.Item(295, 1).Value = 319.2480568410749
$ws.Cells.Item(295, 2).Value = -14.11347504157116
$ws.Cells.Item(296, 1).Value = 322.0598051068358
$ws.Cells.Item(296, 2).Value = -14.2612049869377
$ws.Cells.Item(297, 1).Value = 324.8634982526433
$ws.Cells.Item(297, 2).Value = -14.40932798130446
$ws.Cells.Item(298, 1).Value = 327.6590735852207
$ws.Cells.Item(298, 2).Value = -14.55785369023147
$ws.Cells.Item(299, 1).Value = 330.4464801409294
$ws.Cells.Item(299, 2).Value = -14.70679433081797
$ws.Cells.Item(300, 1).Value = 333.2256316102565
$ws.Cells.Item(300, 2).Value = -14.85615525778927
$ws.Cells.Item(301, 1).Value = 335.9964707867462
$ws.Cells.Item(301, 2).Value = -15.00594791881475
$ws.Cells.Item(302, 1).Value = 338.7589257725257
$ws.Cells.Item(302, 2).Value = -15.15618100081885
$ws.Cells.Item(303, 1).Value = 341.5129330841825
$ws.Cells.Item(303, 2).Value = -15.30686521755972
$ws.Cells.Item(304, 1).Value = 344.2584219814806
$ws.Cells.Item(304, 2).Value = -15.45800952826226
$ws.Cells.Item(305, 1).Value = 346.9953184078549
$ws.Cells.Item(305, 2).Value = -15.60962342668761
$ws.Cells.Item(306, 1).Value = 349.723566078661
$ws.Cells.Item(306, 2).Value = -15.76171934529392
$ws.Cells.Item(307, 1).Value = 352.4430887492692
$ws.Cells.Item(307, 2).Value = -15.91430654240279
$ws.Cells.Item(308, 1).Value = 355.1538265970422
$ws.Cells.Item(308, 2).Value = -16.06739738179427
$ws.Cells.Item(309, 1).Value = 357.8557116656273
$ws.Cells.Item(309, 2).Value = -16.22100351279479
$ws.Cells.Item(310, 1).Value = 360.5486628908087
$ws.Cells.Item(310, 2).Value = -16.37513352694742
$ws.Cells.Item(311, 1).Value = 363.2326305648619
$ws.Cells.Item(311, 2).Value = -16.5298032590882
$ws.Cells.Item(312, 1).Value = 365.907541440833
$ws.Cells.Item(312, 2).Value = -16.68502383816376
$ws.Cells.Item(313, 1).Value = 368.5733191867066
$ws.Cells.Item(313, 2).Value = -16.84080602741139
$ws.Cells.Item(314, 1).Value = 371.2299049985633
$ws.Cells.Item(314, 2).Value = -16.99716470965619
$ws.Cells.Item(315, 1).Value = 373.8772274506714
$ws.Cells.Item(315, 2).Value = -17.1541126506351
$ws.Cells.Item(316, 1).Value = 376.5152191766381
$ws.Cells.Item(316, 2).Value = -17.31166342054392
$ws.Cells.Item(317, 1).Value = 379.143820353212
$ws.Cells.Item(317, 2).Value = -17.4698330606858
$ws.Cells.Item(318, 1).Value = 381.7629531241572
$ws.Cells.Item(318, 2).Value = -17.62863410658922
$ws.Cells.Item(319, 1).Value = 384.3725586182796
$ws.Cells.Item(319, 2).Value = -17.78808319354827
$ws.Cells.Item(320, 1).Value = 386.9725488260881
$ws.Cells.Item(320, 2).Value = -17.94819174659506
$ws.Cells.Item(321, 1).Value = 389.562887656147
$ws.Cells.Item(321, 2).Value = -18.10898265555215
$ws.Cells.Item(322, 1).Value = 392.14348168136
$ws.Cells.Item(322, 2).Value = -18.27046672633634
$ws.Cells.Item(323, 1).Value = 394.7142671881351
$ws.Cells.Item(323, 2).Value = -18.43266183152816
$ws.Cells.Item(324, 1).Value = 397.2751826772408
$ws.Cells.Item(324, 2).Value = -18.5955871319818
$ws.Cells.Item(325, 1).Value = 399.8261513001669
$ws.Cells.Item(325, 2).Value = -18.75925894175774
$ws.Cells.Item(326, 1).Value = 402.3671094265976
$ws.Cells.Item(326, 2).Value = -18.92369655989787
$ws.Cells.Item(327, 1).Value = 404.8979913215935
$ws.Cells.Item(327, 2).Value = -19.08892019577608
$ws.Cells.Item(328, 1).Value = 407.4187255287925
$ws.Cells.Item(328, 2).Value = -19.2549491756184
$ws.Cells.Item(329, 1).Value = 409.9292280051077
$ws.Cells.Item(329, 2).Value = -19.42180064877479
$ws.Cells.Item(330, 1).Value = 412.4294485844933
$ws.Cells.Item(330, 2).Value = -19.58949993000327
$ws.Cells.Item(331, 1).Value = 414.9193143651427
$ws.Cells.Item(331, 2).Value = -19.75806839208979
$ws.Cells.Item(332, 1).Value = 417.3987559370753
$ws.Cells.Item(332, 2).Value = -19.92752866372022
$ws.Cells.Item(333, 1).Value = 419.8676877927811
$ws.Cells.Item(333, 2).Value = -20.0979006144598
$ws.Cells.Item(334, 1).Value = 422.3260546174901
$ws.Cells.Item(334, 2).Value = -20.26921207427822
$ws.Cells.Item(335, 1).Value = 424.773782497233
$ws.Cells.Item(335, 2).Value = -20.44148740139428
$ws.Cells.Item(336, 1).Value = 427.2108005390471
$ws.Cells.Item(336, 2).Value = -20.61475246284351
$ws.Cells.Item(337, 1).Value = 429.6370494468468
$ws.Cells.Item(337, 2).Value = -20.78903710263667
$ws.Cells.Item(338, 1).Value = 432.0524319246359
$ws.Cells.Item(338, 2).Value = -20.96436340913545
$ws.Cells.Item(339, 1).Value = 434.4568947466847
$ws.Cells.Item(339, 2).Value = -21.14076453741118
$ws.Cells.Item(340, 1).Value = 436.8503685072328
$ws.Cells.Item(340, 2).Value = -21.31827115265495
$ws.Cells.Item(341, 1).Value = 439.2327546081752
$ws.Cells.Item(341, 2).Value = -21.4969084225377
$ws.Cells.Item(342, 1).Value = 441.6040199226645
$ws.Cells.Item(342, 2).Value = -21.67671858093513
$ws.Cells.Item(343, 1).Value = 443.9640628531031
$ws.Cells.Item(343, 2).Value = -21.85772821596511
$ws.Cells.Item(344, 1).Value = 446.3128150732211
$ws.Cells.Item(344, 2).Value = -22.03997383733696
$ws.Cells.Item(345, 1).Value = 448.6502083234316
$ws.Cells.Item(345, 2).Value = -22.22349329793606
$ws.Cells.Item(346, 1).Value = 450.9761801456556
$ws.Cells.Item(346, 2).Value = -22.40832766148625
$ws.Cells.Item(347, 1).Value = 453.2906459936422
$ws.Cells.Item(347, 2).Value = -22.59451396968963
$ws.Cells.Item(348, 1).Value = 455.5935294278781
$ws.Cells.Item(348, 2).Value = -22.78209368695408
$ws.Cells.Item(349, 1).Value = 457.8847616684391
$ws.Cells.Item(349, 2).Value = -22.97111106936689
$ws.Cells.Item(350, 1).Value = 460.1642660975166
$ws.Cells.Item(350, 2).Value = -23.16161118360855
$ws.Cells.Item(351, 1).Value = 462.4319596811087
$ws.Cells.Item(351, 2).Value = -23.35363939359294
$ws.Cells.Item(352, 1).Value = 464.6877853317664
$ws.Cells.Item(352, 2).Value = -23.54724958594274
$ws.Cells.Item(353, 1).Value = 466.9316593901643
$ws.Cells.Item(353, 2).Value = -23.74249181204256
$ws.Cells.Item(354, 1).Value = 469.1635083959185
$ws.Cells.Item(354, 2).Value = -23.93942074007701
$ws.Cells.Item(355, 1).Value = 471.3832617120587
$ws.Cells.Item(355, 2).Value = -24.13809477752753
$ws.Cells.Item(356, 1).Value = 473.5908491053529
$ws.Cells.Item(356, 2).Value = -24.33857504961345
$ws.Cells.Item(357, 1).Value = 475.7861775554152
$ws.Cells.Item(357, 2).Value = -24.5409197770988
$ws.Cells.Item(358, 1).Value = 477.969196724789
$ws.Cells.Item(358, 2).Value = -24.74520158261567
$ws.Cells.Item(359, 1).Value = 480.1398239906171
$ws.Cells.Item(359, 2).Value = -24.95148807121053
$ws.Cells.Item(360, 1).Value = 482.2979838085744
$ws.Cells.Item(360, 2).Value = -25.1598522575557
$ws.Cells.Item(361, 1).Value = 484.4436077164808
$ws.Cells.Item(361, 2).Value = -25.37037315200014
$ws.Cells.Item(362, 1).Value = 486.5766382787576
$ws.Cells.Item(362, 2).Value = -25.58313677494359
$ws.Cells.Item(363, 1).Value = 488.6969859144046
$ws.Cells.Item(363, 2).Value = -25.79822507911286
$ws.Cells.Item(364, 1).Value = 490.8045891470704
$ws.Cells.Item(364, 2).Value = -26.01573196789037
$ws.Cells.Item(365, 1).Value = 492.8993847100127
$ws.Cells.Item(365, 2).Value = -26.23575617965308
$ws.Cells.Item(366, 1).Value = 494.9813209202326
$ws.Cells.Item(366, 2).Value = -26.45840618934429
$ws.Cells.Item(367, 1).Value = 497.0503225650972
$ws.Cells.Item(367, 2).Value = -26.68378830673663
$ws.Cells.Item(368, 1).Value = 499.1063266070253
$ws.Cells.Item(368, 2).Value = -26.912019792966
$ws.Cells.Item(369, 1).Value = 501.1492972948251
$ws.Cells.Item(369, 2).Value = -27.14323189337003
$ws.Cells.Item(370, 1).Value = 503.1791588015211
$ws.Cells.Item(370, 2).Value = -27.37755240497425
$ws.Cells.Item(371, 1).Value = 505.1958805930281
$ws.Cells.Item(371, 2).Value = -27.61512990553885
$ws.Cells.Item(372, 1).Value = 507.1994241214705
$ws.Cells.Item(372, 2).Value = -27.85612020659054
$ws.Cells.Item(373, 1).Value = 509.1897369028907
$ws.Cells.Item(373, 2).Value = -28.10068338567135
$ws.Cells.Item(374, 1).Value = 511.166825332004
$ws.Cells.Item(374, 2).Value = -28.34900920466705
$ws.Cells.Item(375, 1).Value = 513.1306502088964
$ws.Cells.Item(375, 2).Value = -28.60128444551788
$ws.Cells.Item(376, 1).Value = 515.0812030706855
$ws.Cells.Item(376, 2).Value = -28.85771670120796
$ws.Cells.Item(377, 1).Value = 517.0185057334626
$ws.Cells.Item(377, 2).Value = -29.1185387784277
$ws.Cells.Item(378, 1).Value = 518.9425738448315
$ws.Cells.Item(378, 2).Value = -29.38399565805013
$ws.Cells.Item(379, 1).Value = 520.8534295114689
$ws.Cells.Item(379, 2).Value = -29.65435050735103
$ws.Cells.Item(380, 1).Value = 522.7511582788068
$ws.Cells.Item(380, 2).Value = -29.92990701010326
$ws.Cells.Item(381, 1).Value = 524.6358312734219
$ws.Cells.Item(381, 2).Value = -30.21098389307208
$ws.Cells.Item(382, 1).Value = 526.5075478920714
$ws.Cells.Item(382, 2).Value = -30.4979313031858
